# "performed lstm and bilstm" - add Embeddings (FastText) and Scaling (MinMax)
# columns/results to the results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (G = Embeddings, H = Scaling).
# Set in this order so new shared strings are interned as:
#   8=Embeddings, 9=FastText, 10=Scaling, 11=MinMax (matches source order)
$ws.Range("G1").Value = "Embeddings"
$ws.Range("G2").Value = "FastText"
$ws.Range("H1").Value = "Scaling"
$ws.Range("H2").Value = "MinMax"

# row, Precision(B), Recall(C), F1 Score(D), Accuracy(E)
$data = @(
    @(2, 0,     0,    0,    0.46),
    @(3, 0,     0,    0,    0.46),
    @(4, 0,     0,    0,    0.44),
    @(5, 0.65,  0.09, 0.16, 0.48),
    @(6, 0.569, 0.9,  0.698,0.58),
    @(7, 0,     0,    0,    0.46),
    @(8, 0.57,  0.76, 0.65, 0.56),
    @(9, 0,     0,    0,    0.45)
)

foreach ($rowData in $data) {
    $r = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
    $ws.Cells.Item($r, 7).Value = "FastText"
    $ws.Cells.Item($r, 8).Value = "MinMax"
}

# Leave selection where the author's cursor ended up after editing
$ws.Range("D19").Select()
